# Generate Report for Handoff
#
# Adds two newly-handed-off files (3aaecde3-...md and f87b0403-...md) to the
# localization status workbook: one new row per file on the "Overview"
# sheet, and one new row per file on each of the "zh-cn" / "de-de" detail
# sheets.

$wb = $excel.ActiveWorkbook

$mdGuid1 = "3aaecde3-bdf9-43d1-9b3a-ef458987af9e"
$mdGuid2 = "f87b0403-a2fb-45fc-8896-bb6dc5a7ca24"

$xlfHash1 = "88607df0c396f6418c68da0d26cbc02fb3194cf1"
$xlfHash2 = "3509956661bc83277bb1ad56e6842879de865853"

$handoffDateTimeZhCn = "2016-03-22 10:40:57"
$handoffDateTimeDeDe = "2016-03-22 10:41:00"
$overviewHandoffDate = "2016-03-22 10:41:00"
$noHandbackDateTime = "0001-01-01 00:00:00"

$status = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$rows = @(
    @{ Row = 5; Guid = $mdGuid1 },
    @{ Row = 6; Guid = $mdGuid2 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $guid = $r.Guid
    $mdName = $guid + ".md"

    $wsOverview.Cells.Item($row, 1).Value = $mdName
    $wsOverview.Hyperlinks.Add(
        $wsOverview.Cells.Item($row, 1),
        "https://github.com/OpenLocalizationTest/oltest/blob/f8247e2ae60c5b7c9a54f28be8899123793559c9/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $wsOverview.Cells.Item($row, 2).Value = $status
    $wsOverview.Cells.Item($row, 3).Value = $status

    $wsOverview.Cells.Item($row, 4).Value = $overviewHandoffDate
    $wsOverview.Cells.Item($row, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnRows = @(
    @{ Row = 5; Guid = $mdGuid1; XlfHash = $xlfHash1 },
    @{ Row = 6; Guid = $mdGuid2; XlfHash = $xlfHash2 }
)

foreach ($r in $zhCnRows) {
    $row = $r.Row
    $guid = $r.Guid
    $mdName = $guid + ".md"
    $xlfName = $guid + "." + $r.XlfHash + ".zh-cn.xlf"

    $wsZhCn.Cells.Item($row, 1).Value = $mdName
    $wsZhCn.Hyperlinks.Add(
        $wsZhCn.Cells.Item($row, 1),
        "https://github.com/OpenLocalizationTest/oltest/blob/f8247e2ae60c5b7c9a54f28be8899123793559c9/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $wsZhCn.Cells.Item($row, 2).Value = ".md"
    $wsZhCn.Cells.Item($row, 3).Value = $status

    $wsZhCn.Cells.Item($row, 4).Value = $xlfName
    $wsZhCn.Hyperlinks.Add(
        $wsZhCn.Cells.Item($row, 4),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfName",
        "",
        "",
        $xlfName
    ) | Out-Null

    $wsZhCn.Cells.Item($row, 5).Value = $handoffDateTimeZhCn
    $wsZhCn.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $wsZhCn.Cells.Item($row, 8).Value = $noHandbackDateTime
    $wsZhCn.Cells.Item($row, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $wsZhCn.Cells.Item($row, 10).Value = "Include"
}

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeRows = @(
    @{ Row = 5; Guid = $mdGuid1; XlfHash = $xlfHash1 },
    @{ Row = 6; Guid = $mdGuid2; XlfHash = $xlfHash2 }
)

foreach ($r in $deDeRows) {
    $row = $r.Row
    $guid = $r.Guid
    $mdName = $guid + ".md"
    $xlfName = $guid + "." + $r.XlfHash + ".de-de.xlf"

    $wsDeDe.Cells.Item($row, 1).Value = $mdName
    $wsDeDe.Hyperlinks.Add(
        $wsDeDe.Cells.Item($row, 1),
        "https://github.com/OpenLocalizationTest/oltest/blob/f8247e2ae60c5b7c9a54f28be8899123793559c9/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $wsDeDe.Cells.Item($row, 2).Value = ".md"
    $wsDeDe.Cells.Item($row, 3).Value = $status

    $wsDeDe.Cells.Item($row, 4).Value = $xlfName
    $wsDeDe.Hyperlinks.Add(
        $wsDeDe.Cells.Item($row, 4),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfName",
        "",
        "",
        $xlfName
    ) | Out-Null

    $wsDeDe.Cells.Item($row, 5).Value = $handoffDateTimeDeDe
    $wsDeDe.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $wsDeDe.Cells.Item($row, 8).Value = $noHandbackDateTime
    $wsDeDe.Cells.Item($row, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $wsDeDe.Cells.Item($row, 10).Value = "Include"
}

"Report rows added for handoff."
